$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

# Update PLC data 2025-10-13 13:59:17
$ws.Range("C2").Value = 273
$ws.Range("C3").Value = 170415
$ws.Range("C4").Value = 161238
$ws.Range("C8").Value = 65.62
